$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.348.07"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.687.24"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5463"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2726"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06447"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07676"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.693.79"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.536"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5812"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008346"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "26.389.33"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.946"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.225"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.011"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.64"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1317"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.895"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06351"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.413"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.95%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.579"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.574"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.675"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6156"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.410"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.238"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Value = "1.115.16"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01626"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8795"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.05"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.839.70"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.199"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05269"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4305"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.039"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.43%  "
